# Updated cryptos list on Sun Nov 17 19:52:06 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: several "Price" values are plain-looking decimal numbers
# (e.g. "235.47", "0.360"). Assigning such a string straight to .Value
# would make Excel auto-detect it as a Number (and normalize e.g.
# "0.360" -> 0.36, "9.10" -> 9.1), but the source data keeps these as
# literal text. To preserve the exact text we briefly mark the cell as
# Text ("@") before writing the value, then restore the "Normal" style
# so no stray number-format style is left behind on the cell.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "90.155.82"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.096.06"
$ws.Range("E3").Value = "  -2.30%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.00%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.26%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -12.32%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.360"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.38%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.093.02"
$ws.Range("E10").Value = "  -2.34%  "

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.719"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.13%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -3.13%  "

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.40%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.52%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "89.911.06"
$ws.Range("E15").Value = "  -0.99%  "

# Row 16 - Toncoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.68%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.653.77"
$ws.Range("E17").Value = "  -2.77%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.086.57"
$ws.Range("E18").Value = "  -3.65%  "

# Row 19 - SuiNetwork
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.43%  "

# Row 20 - PEPE
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000212"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.11%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "433.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.73%  "

# Row 23 - Polkadot
$ws.Range("E23").Value = "  +3.71%  "

# Row 24 - Uniswap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.75%  "

# Row 25 - NEARProtocol
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.12%  "

# Row 26 - Litecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.49%  "

# Row 27 - Aptos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.54%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  +2.85%  "

# Row 32 - Cronos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.157"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.15%  "

# Row 33 - Stellar
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.193"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.39%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.26%  "

# Row 35 - Kaspa
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.151"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.36%  "

# Row 36 - dogwifhat
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.79%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.48%  "

# Row 38 - Bittensor
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "498.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.13%  "

# Row 39 - PancakeSwap
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.84%  "

# Row 40 - Fetch.AI
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.16%  "

# Row 41 - MantraDAO
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +59.23%  "

# Row 42 - Hedera
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0867"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.62%  "

# Row 43 - WhiteBITCoin
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "

# Row 45 - PolygonEcosystemToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.400"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.38%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  -6.06%  "

# Row 47 - now Monero (was ARBITRUM)
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "

# Row 48 - now ARBITRUM (was Monero)
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.681"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.57%  "

# Row 49 - OKB
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.37%  "

# Row 50 - ImmutableX
$ws.Range("E50").Value = "  -4.54%  "

# Row 51 - now Filecoin (was FirstDigitalUSD)
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.99%  "
